# Generate Report for Handback
#
# Applies the "handback" localization-status update:
#  - Overview + per-language "Status" text changes from "In Translation" to
#    "Handed back: in sync with en-US"
#  - per-language sheets (zh-cn / de-de) get their "Latest Target File" /
#    "Latest Handback File" / "Latest Handback DateTime" columns populated
#    for the two rows, including a hyperlink on the new "Latest Target File"
#    cell pointing at the same source-markdown URL as column A's hyperlink.
#  - column widths widen to fit the newly-populated / newly-lengthened text.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$srcUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2d3c1edbb5630d861cb70aa07bc3d02c4581faa3/e2e/ccc8df8a-9572-4aa2-81c5-f987c5086792.md"
$srcUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2d3c1edbb5630d861cb70aa07bc3d02c4581faa3/e2e/cdb5e241-d7fc-45a6-bca2-fb9b5aed190b.md"
$srcName1 = "ccc8df8a-9572-4aa2-81c5-f987c5086792.md"
$srcName2 = "cdb5e241-d7fc-45a6-bca2-fb9b5aed190b.md"

# ---------------------------------------------------------------------
# Overview sheet: status text for both rows/languages, then widen the
# zh-cn / de-de columns (E, F) to fit it.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = 29.144371396019366
$wsOverview.Columns.Item(6).ColumnWidth = 29.144371396019366

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsZh.Range("I2").Value = $srcName1
$wsZh.Range("J2").Value = "ccc8df8a-9572-4aa2-81c5-f987c5086792.66445827c868bcf48a3f0dd029aca828b6dba094.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-17 06:20:40"

$wsZh.Range("I3").Value = $srcName2
$wsZh.Range("J3").Value = "cdb5e241-d7fc-45a6-bca2-fb9b5aed190b.5061da1696e3becaddd1b4387fc0684f512aa181.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-17 06:20:40"

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $srcUrl1, "", "", $srcName1)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $srcUrl2, "", "", $srcName2)

$wsZh.Columns.Item(3).ColumnWidth = 29.144371396019366
$wsZh.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsZh.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

$wsDe.Range("I2").Value = $srcName1
$wsDe.Range("J2").Value = "ccc8df8a-9572-4aa2-81c5-f987c5086792.66445827c868bcf48a3f0dd029aca828b6dba094.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-17 06:20:47"

$wsDe.Range("I3").Value = $srcName2
$wsDe.Range("J3").Value = "cdb5e241-d7fc-45a6-bca2-fb9b5aed190b.5061da1696e3becaddd1b4387fc0684f512aa181.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-17 06:20:47"

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $srcUrl1, "", "", $srcName1)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $srcUrl2, "", "", $srcName2)

$wsDe.Columns.Item(3).ColumnWidth = 29.144371396019366
$wsDe.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsDe.Columns.Item(10).ColumnWidth = 39.166666666666664
